$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 21,9
$arr[0,0] = 0.5089441903256321
$arr[0,1] = 0.4943333713621592
$arr[0,2] = 0
$arr[0,3] = 0.3246404326649622
$arr[0,4] = 0.3586383285908433
$arr[0,5] = 0.3586383285908433
$arr[0,6] = 0.1794793107198907
$arr[0,7] = 0.2610012569162791
$arr[0,8] = 0.5942751998756042
$arr[1,0] = 0.3923116926485909
$arr[1,1] = 0
$arr[1,2] = 0
$arr[1,3] = 0.4908480066013919
$arr[1,4] = 0.554957660978823
$arr[1,5] = 0.554957660978823
$arr[1,6] = 0.6500194412739141
$arr[1,7] = 0.6102808058241208
$arr[1,8] = 0.4686367749667203
$arr[2,0] = 0.224450920129413
$arr[2,1] = 0.5179264704491768
$arr[2,2] = 0.05103467618520385
$arr[2,3] = 0.2655360114382984
$arr[2,4] = 0.1887897996101501
$arr[2,5] = 0.1887897996101501
$arr[2,6] = 0.2136556847703774
$arr[2,7] = 0.1893997023084638
$arr[2,8] = 0.09180260090716087
$arr[3,0] = 0.7780876219875478
$arr[3,1] = -0.02472840254945134
$arr[3,2] = 0
$arr[3,3] = 0.8490787306869069
$arr[3,4] = 0.7437040793122646
$arr[3,5] = 0.7437040793122646
$arr[3,6] = 0.8490787306869069
$arr[3,7] = 0.7780876219875478
$arr[3,8] = 0.8490787306869069
$arr[4,0] = 0.2220561864913329
$arr[4,1] = 0.1612076331529478
$arr[4,2] = 0.1952695269526953
$arr[4,3] = 0.1442483262325015
$arr[4,4] = 0.2865013774104683
$arr[4,5] = 0.2865013774104683
$arr[4,6] = 0.2102960478343299
$arr[4,7] = 0.3163329161451815
$arr[4,8] = 0.2907662082514735
$arr[5,0] = 0.08588162782127758
$arr[5,1] = 0.07482671326787527
$arr[5,2] = 0.02355523043160223
$arr[5,3] = 0.2511841651758604
$arr[5,4] = 0.03270642201834869
$arr[5,5] = 0.03270642201834869
$arr[5,6] = -0.02103211009174305
$arr[5,7] = -0.02615503290715429
$arr[5,8] = 0.2619408900464274
$arr[6,0] = 0.1726881347951491
$arr[6,1] = 0.1923971217807718
$arr[6,2] = 0
$arr[6,3] = 0.2965919384057971
$arr[6,4] = -0.05504256060165497
$arr[6,5] = -0.05504256060165497
$arr[6,6] = 0.4315771524583482
$arr[6,7] = 0.4495362037934376
$arr[6,8] = 0.2145429743060812
$arr[7,0] = 0.01517832801143781
$arr[7,1] = 0
$arr[7,2] = 0
$arr[7,3] = 0.8228271770322414
$arr[7,4] = 0.6003473083602083
$arr[7,5] = 0.6003473083602083
$arr[7,6] = 0.1305616848814452
$arr[7,7] = 0.9200414896919314
$arr[7,8] = 0.06494837344412167
$arr[8,0] = 0.04248245289988915
$arr[8,1] = 0
$arr[8,2] = 0
$arr[8,3] = 0.09954158480681068
$arr[8,4] = -0.05762933857236419
$arr[8,5] = -0.05762933857236419
$arr[8,6] = 0.06033940917661845
$arr[8,7] = -0.00678179402922467
$arr[8,8] = 0.140262599260948
$arr[9,0] = 0.08434343434343435
$arr[9,1] = 0
$arr[9,2] = -0.004700507972343765
$arr[9,3] = 0.1019917402896127
$arr[9,4] = 0.04508496450849644
$arr[9,5] = 0.04508496450849644
$arr[9,6] = 0.04427525712697285
$arr[9,7] = 0.1069334692837114
$arr[9,8] = 0.09799751395704749
$arr[10,0] = -0.03583473861720073
$arr[10,1] = 0
$arr[10,2] = -0.005319148936170364
$arr[10,3] = 0.002099958000840043
$arr[10,4] = -0.03409090909090909
$arr[10,5] = -0.03409090909090909
$arr[10,6] = 0.01451679800912478
$arr[10,7] = -0.001686340640809492
$arr[10,8] = -0.02574926129168418
$arr[11,0] = -0.08450704225352125
$arr[11,1] = 0
$arr[11,2] = 0
$arr[11,3] = 0.05688622754491006
$arr[11,4] = -0.024390243902439
$arr[11,5] = -0.024390243902439
$arr[11,6] = -0.0405405405405406
$arr[11,7] = 0.1812865497076024
$arr[11,8] = 0.01355013550135503
$arr[12,0] = 0.04281767955801098
$arr[12,1] = 0
$arr[12,2] = 0.03076923076923077
$arr[12,3] = -0.003429355281207231
$arr[12,4] = 0.1970802919708028
$arr[12,5] = 0.1970802919708028
$arr[12,6] = -0.01167883211678843
$arr[12,7] = 0.04281767955801098
$arr[12,8] = 0
$arr[13,0] = 0.3743869209809265
$arr[13,1] = 0.0143742255266419
$arr[13,2] = 0.165769000598444
$arr[13,3] = 0.2144776766363242
$arr[13,4] = 0.3024827024827025
$arr[13,5] = 0.3024827024827025
$arr[13,6] = 0.5777987718469533
$arr[13,7] = 0.1857914416285833
$arr[13,8] = 0.1987676056338028
$arr[14,0] = 0.3624811388968273
$arr[14,1] = 0.03490904536133276
$arr[14,2] = -0.001629632488204126
$arr[14,3] = 0.4929564166194816
$arr[14,4] = 0.2456420626559125
$arr[14,5] = 0.2456420626559125
$arr[14,6] = 0.3528218428011595
$arr[14,7] = 0.3426105006068903
$arr[14,8] = 0.2234327044307217
$arr[15,0] = 0.1467170656586868
$arr[15,1] = 0
$arr[15,2] = 0.1025292357900461
$arr[15,3] = 0.01357658039881208
$arr[15,4] = -0.03959683225341973
$arr[15,5] = -0.03959683225341973
$arr[15,6] = 0.1421051149959182
$arr[15,7] = -0.006508991947928046
$arr[15,8] = 0.1447084233261339
$arr[16,0] = 0.03518437202263319
$arr[16,1] = 0.03138581299500854
$arr[16,2] = 0.005357593271323863
$arr[16,3] = 0.0293527236723549
$arr[16,4] = 0.03489933096170147
$arr[16,5] = 0.03489933096170147
$arr[16,6] = 0.0298199743222094
$arr[16,7] = 0.06446305198588806
$arr[16,8] = 0.06483324433984185
$arr[17,0] = 0.4246698766347226
$arr[17,1] = -0.001682722548086535
$arr[17,2] = 0
$arr[17,3] = 0.4314287805532591
$arr[17,4] = 0.3036008406515629
$arr[17,5] = 0.3036008406515629
$arr[17,6] = 0.3907868125394266
$arr[17,7] = 0.446825073863188
$arr[17,8] = 0.4493924241640341
$arr[18,0] = -0.02035278154681141
$arr[18,1] = 0
$arr[18,2] = 0
$arr[18,3] = 0.4327097163548582
$arr[18,4] = 0.3268505955073119
$arr[18,5] = 0.3268505955073119
$arr[18,6] = -0.01992678416424493
$arr[18,7] = 0.04226415094339623
$arr[18,8] = -0.01403138201569091
$arr[19,0] = -0.05106795963914094
$arr[19,1] = 0
$arr[19,2] = 0
$arr[19,3] = -0.01851905628050766
$arr[19,4] = 0.1062404390290867
$arr[19,5] = 0.1062404390290867
$arr[19,6] = 0.0005100182149363802
$arr[19,7] = 0.1394124820813454
$arr[19,8] = 0.1394124820813454
$arr[20,0] = -0.0145470234833475
$arr[20,1] = -0.005530368669925899
$arr[20,2] = -0.007376670918393435
$arr[20,3] = 0.002157865200333181
$arr[20,4] = 0.001800914066815311
$arr[20,5] = 0.001800914066815311
$arr[20,6] = -0.007235863445365188
$arr[20,7] = -0.007917137795485647
$arr[20,8] = 0.02098811055254978

$ws.Range("B2:J22").Value = $arr
Write-Output "applied grid"